# Add "TEST_SHEET", "TEXEL_SHEET_TRACKER" and "TEST_SCALAR_INPUT" sheets
# around the existing "Sheet1", matching a workbook that tracks its own
# sheet inventory and carries a couple of scalar test inputs.
#
# NOTE: worksheet object handles go stale (silently re-target whatever is
# the *current* ActiveSheet) as soon as another Worksheets.Add()/Delete()
# call happens anywhere in the workbook. So: (1) do every sheet add/remove
# first, with no interleaved cell writes, then (2) re-fetch each sheet via
# Worksheets.Item("<name>") right before touching it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Phase 1: create the sheets, in an order that reproduces the target
# sheetIds (sheetId = current-max-alive-id + 1 at the moment of Add, so
# three throw-away sheets push the counter up to land TEST_SCALAR_INPUT
# on sheetId 9).
# ---------------------------------------------------------------------
$tracker = $wb.Worksheets.Add($wb.Worksheets.Item("Sheet1"))
$tracker.Name = "TEXEL_SHEET_TRACKER"

$testSheet = $wb.Worksheets.Add($wb.Worksheets.Item("TEXEL_SHEET_TRACKER"))
$testSheet.Name = "TEST_SHEET"

$spacer1 = $wb.Worksheets.Add()
$spacer1Name = $spacer1.Name
$spacer2 = $wb.Worksheets.Add()
$spacer2Name = $spacer2.Name
$spacer3 = $wb.Worksheets.Add()
$spacer3Name = $spacer3.Name

$scalarSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item("Sheet1"))
$scalarSheet.Name = "TEST_SCALAR_INPUT"

$wb.Worksheets.Item($spacer1Name).Delete()
$wb.Worksheets.Item($spacer2Name).Delete()
$wb.Worksheets.Item($spacer3Name).Delete()

# ---------------------------------------------------------------------
# Phase 2: populate each sheet. No more Add()/Delete() calls from here
# on, so re-fetching by name keeps every reference fresh.
# ---------------------------------------------------------------------

# --- TEST_SHEET: blank placeholder sheet, just move the selection ---
$wb.Worksheets.Item("TEST_SHEET").Activate()
$wb.Worksheets.Item("TEST_SHEET").Range("D17").Select()

# --- TEXEL_SHEET_TRACKER: header row + one data row describing TEST_SHEET ---
$wb.Worksheets.Item("TEXEL_SHEET_TRACKER").Activate()
$excel.ActiveWindow.DisplayGridlines = $false

$wb.Worksheets.Item("TEXEL_SHEET_TRACKER").Range("A1").Value = "sheet_name"
$wb.Worksheets.Item("TEXEL_SHEET_TRACKER").Range("B1").Value = "descr"
$wb.Worksheets.Item("TEXEL_SHEET_TRACKER").Range("A2").Value = "TEST_SHEET"
$wb.Worksheets.Item("TEXEL_SHEET_TRACKER").Range("B2").Value = "first sheet added"

$helper = $wb.Worksheets.Item("TEXEL_SHEET_TRACKER").Range("Z100")
$helper.Font.Bold = $true
$helper.Interior.Color = 12566463
$helper.HorizontalAlignment = -4108
$helper.Copy()
$wb.Worksheets.Item("TEXEL_SHEET_TRACKER").Range("A1:B1").PasteSpecial(-4122)
$helper.Clear()
$excel.CutCopyMode = $false

$wb.Worksheets.Item("TEXEL_SHEET_TRACKER").Columns.Item(1).EntireColumn.AutoFit()
$wb.Worksheets.Item("TEXEL_SHEET_TRACKER").Columns.Item(2).EntireColumn.AutoFit()

$wb.Worksheets.Item("TEXEL_SHEET_TRACKER").PageSetup.Orientation = 1

$wb.Worksheets.Item("TEXEL_SHEET_TRACKER").Range("B2").Select()

# --- Sheet1: untouched content, selection moved ---
$wb.Worksheets.Item("Sheet1").Activate()
$wb.Worksheets.Item("Sheet1").Range("F20").Select()

# --- TEST_SCALAR_INPUT: var_name/value table, becomes the active sheet ---
$wb.Worksheets.Item("TEST_SCALAR_INPUT").Activate()

$wb.Worksheets.Item("TEST_SCALAR_INPUT").Range("A1").Value = "var_name"
$wb.Worksheets.Item("TEST_SCALAR_INPUT").Range("B1").Value = "value"
$wb.Worksheets.Item("TEST_SCALAR_INPUT").Range("A2").Value = "a"
$wb.Worksheets.Item("TEST_SCALAR_INPUT").Range("B2").Value = 1
$wb.Worksheets.Item("TEST_SCALAR_INPUT").Range("A3").Value = "b"
$wb.Worksheets.Item("TEST_SCALAR_INPUT").Range("B3").Value = 2

$wb.Worksheets.Item("TEST_SCALAR_INPUT").Range("B4").Select()
